$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph: "License Information" heading -> plain "Aquifer Open
# Study Notes (Book Intros)" bold run, no heading style.
# ------------------------------------------------------------------
$pHeading = $d.Paragraphs.Item(4)
$pHeading.Range.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "Aquifer Open Study Notes (Book Intros)", 2)
$d.Paragraphs.Item(4).Style = "Normal"
$d.Paragraphs.Item(4).Range.Bold = 1

# ------------------------------------------------------------------
# Paragraph: main license / attribution paragraph.
# ------------------------------------------------------------------
$pLic = $d.Paragraphs.Item(5)

# Drop the old lead-in ("<bold title> (French) is based on: ") so the
# "Tyndale Open Study Notes" run survives untouched.
$pLic.Range.Find.Execute("Notes d'étude - Introductions aux livres (Tyndale) (French) is based on: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Insert the new introductory sentence in front of "Tyndale Open Study Notes".
$pLicA = $d.Paragraphs.Item(5)
$insA = $pLicA.Range.Duplicate
$insA.Collapse(1)
$insA.Text = "This work is an adaptation of "

# Replace the old ", Tyndale House Publishers, 2019, ... license." tail
# (which spans both hyperlinks) with the new connector sentence.
$pLicB = $d.Paragraphs.Item(5)
$pLicB.Range.Find.Execute(", Tyndale House Publishers, 2019, which is licensed under a CC BY-SA 4.0 license.", $true, $false, $false, $false, $false, $true, 1, $false, " © 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, ", 2)

# Append the new closing sentence (with its own "Aquifer Open Study Notes" mention).
$pLicC = $d.Paragraphs.Item(5)
$endC = $pLicC.Range.Duplicate
$endC.MoveEnd(1, -1)
$endC.Collapse(0)
$endC.Text = "Aquifer Open Study Notes, was created by Mission Mutual and is also licensed under CC BY-SA 4.0."

# ------------------------------------------------------------------
# Paragraph: "This PDF version is provided under the same license."
# -> new multi-language adaptation blurb.
# ------------------------------------------------------------------
$pPdf = $d.Paragraphs.Item(6)
$pPdf.Range.Find.Execute("This PDF version is provided under the same license.", $true, $false, $false, $false, $false, $true, 1, $false, "This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文).", 2)
